$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Copy formatting from the last existing row (14) down to the new row (15),
# one column at a time so each cell keeps the style used by that column.
$cols = @("A", "B", "C", "D", "E")
foreach ($col in $cols) {
    $ws.Range($col + "14").Copy()
    $ws.Range($col + "15").PasteSpecial(-4122)  # xlPasteFormats
}

# New withdrawal record dated 8 May 2018.
$ws.Range("A15").Value = [DateTime]::ParseExact("2018-05-08", "yyyy-MM-dd", $null)
$ws.Range("B15").Value = "Marsicovetere Maria"
$ws.Range("C15").Value = "Tessuto a quadri"
$ws.Range("D15").Value = "Mt."
$ws.Range("E15").Value = 1.1
